$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "FindDaggerRoom" entry from the map (cell F7)
$ws.Range("F7").Clear()

# Update the active selection to match the author's final cursor position
$ws.Range("F6").Select()
